$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Rows 2-14: "VBP" -> "Valor bruto da produção industrial"
$ws.Range("B2:B14").Value = "Valor bruto da produção industrial"

# Rows 15-27: "Custo das Operações" -> "Custo das operações industriais"
$ws.Range("B15:B27").Value = "Custo das operações industriais"

# Rows 28-40: "Valor da Transformação" -> "Valor da transformação industrial"
$ws.Range("B28:B40").Value = "Valor da transformação industrial"
